$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Compound" column (E) for rows 21-28 from a blank space to "FS".
# The "Model_Name" column (J) is a shared formula referencing column E, so it
# will recompute automatically.
for ($r = 21; $r -le 28; $r++) {
    $ws.Cells.Item($r, 5).Value = "FS"
}

# Update the view state to match the saved selection/scroll position.
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E28").Select()
